# The "Test Orientation" column (F) — header + all per-specimen "Blue on
# Top"/"Blue on Bottom" values — is being retired from this coupon sheet,
# so its contents (F1:F21) are cleared out entirely. Columns G onward
# (Test File Name, Run number headers / data, etc.) are left exactly
# where they are.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:F21").ClearContents() | Out-Null

# Move the active selection to reflect where the author was last working.
$ws.Range("F9").Select() | Out-Null
